$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and the Aave/TrustWalletToken row swap)
# Cells whose new value could be misread as a number (e.g. "1.00", "11.10") are
# written via a text-format round-trip so Excel keeps them as literal text, matching
# the original inline-string cell type, then ClearFormats() drops the temporary
# number-format so no stray style index is left on the cell.
$ws.Range("D2").Value = "34.899.63"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.811.29"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.35"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.315"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "2.071.32"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.810.94"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.69"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.657"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("D17").Value = "34.847.59"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").Value = "0.0₃0784"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +6.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.64"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.58%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  +6.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.69"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  +31.44%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "3.339.15"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0552"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.88"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("B36").Value = "Aave"
$ws.Range("C36").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "93.29"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.29%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.13"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.678"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("D40").Value = "1.306.95"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.28"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.986"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.80"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("E47").Value = "  +7.90%  "
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D49").Value = "1.986.00"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  +5.58%  "
